$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reorder the attribute/type rows (rows 3-21; row 2 "case"/"str" is unchanged) ---
$rows = @(
    @("case", "str"),
    @("planned_operation_time", "str"),
    @("time:timestamp", "datetime"),
    @("lifecycle:transition", "str"),
    @("org:resource", "str"),
    @("case:concept:name", "str"),
    @("response_status_code", "float"),
    @("current_task", "str"),
    @("requested_service_url", "str"),
    @("event_id", "str"),
    @("parameters", "dict"),
    @("process_model_id", "str"),
    @("SubProcessID", "str"),
    @("concept:name", "str"),
    @("identifier:id", "str"),
    @("lifecycle:state", "str"),
    @("operation_end_time", "datetime"),
    @("unsatisfied_condition_description", "str"),
    @("human_workstation_green_button_pressed", "float"),
    @("complete_service_time", "str")
)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# --- Strip the bold/bordered/centered header style from A1:B1 so they fall back
#     to the workbook's default (unstyled) cell format. ---
$ws.Cells.ClearFormats()
